# Actualizacion automatica 2025-09-12 15:05:09
#
# A new advisor/client, "SALAZAR VERA ENRIQUE WILLIAM", is inserted (in its
# alphabetically-sorted position) between "QUIJIJE MENDOZA GENESIS XIOMARA"
# and "SOLIS OCAMPO DIMAS ABDON" on both the "VENTAS POR GRUPO" and
# "VENTA MENSUAL" sheets. That pushes every following data row down by one
# and grows the trailing summary/count row from row 42 to row 43 (and the
# "X de 40" counters on "VENTAS POR GRUPO" become "X de 41").

$wb = $excel.ActiveWorkbook

$newName = "SALAZAR VERA ENRIQUE WILLIAM"
$officeName = "OFICINA-CATAECSA"
$insertRow = 36

# ----- Sheet 1: VENTAS POR GRUPO (columns A:R) -----
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Push rows 36..41 down to 37..42 (and the old summary row 42 becomes 43),
# inheriting formatting from the row being displaced.
$ws1.Rows.Item($insertRow).Insert()

$ws1.Cells.Item($insertRow, 1).Value = $officeName
$ws1.Cells.Item($insertRow, 2).Value = $newName
for ($c = 3; $c -le 18; $c++) {
    $ws1.Cells.Item($insertRow, $c).Value = 0
}

# The summary row (previously r=42, now r=43) counts "X de 40" -> "X de 41"
# since the roster grew by one advisor.
$summaryRow1 = 43
for ($c = 3; $c -le 18; $c++) {
    $cell = $ws1.Cells.Item($summaryRow1, $c)
    $cell.Value = $cell.Text -replace "de 40", "de 41"
}

# ----- Sheet 2: VENTA MENSUAL (columns A:G) -----
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item($insertRow).Insert()

$ws2.Cells.Item($insertRow, 1).Value = $officeName
$ws2.Cells.Item($insertRow, 2).Value = $newName
for ($c = 3; $c -le 7; $c++) {
    $ws2.Cells.Item($insertRow, $c).Value = 0
}
